$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3)
$ws.Cells.Item(3, 2).Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (row 8)
$ws.Cells.Item(8, 2).Value = "2025-10-29T22:15:57+01:00"

# Insert a new "Jurisdiction" property row right after "Contact" (row 10),
# pushing Description and everything below it down by one row.
$ws.Rows.Item(11).Insert()

# Match the formatting (borders/style) used by the surrounding property rows.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""
